# Commit: "added tests, few comments and behavior"
#
# The underlying model only has 9 true input cells on Sheet1 (everything
# else is a formula that recalculates automatically). This script updates
# those 9 inputs; Excel's recalculation engine takes care of every
# dependent formula cell (F2, J2, N2, D3, F3, H3, J3, L3, D16, F16, ...
# all the way down through the val/grad cascade in rows 16-41).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- row 2 block -----------------------------------------------------
$ws.Range("B2").Value  = 2.2000000000000002
$ws.Range("D2").Value  = 2.5
$ws.Range("H2").Value  = 2.5
$ws.Range("L2").Value  = 2.5

# --- row 4 ------------------------------------------------------------
$ws.Range("H4").Value  = 2.5

# --- row 7 block -------------------------------------------------------
$ws.Range("H7").Value  = 2.5
$ws.Range("L7").Value  = 2.5

# --- row 9 --------------------------------------------------------------
$ws.Range("H9").Value  = 2.5

# --- row 12 ---------------------------------------------------------------
$ws.Range("L12").Value = 2.5

# Match the author's last on-screen selection (cell L13) as recorded in
# the saved sheetView/selection of the worksheet XML.
$ws.Range("L13").Select()
